$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 349 - this shifts the existing rows
# 349-378 down to 351-380, carrying their values/styles with them.
$ws.Rows.Item(349).Resize(2).Insert()

# Row 349: new weekly observation, "Primera" quality (same categorical
# fields as the row that used to sit here, now at 351).
$ws.Cells.Item(349, 1).Value = 6
$ws.Cells.Item(349, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(349, 3).Value = "Metropolitana"
$ws.Cells.Item(349, 4).Value = 44461
$ws.Cells.Item(349, 5).Value = 13
$ws.Cells.Item(349, 6).Value = 100112017
$ws.Cells.Item(349, 7).Value = "Apio"
$ws.Cells.Item(349, 8).Value = "Americana (o)"
$ws.Cells.Item(349, 9).Value = "Primera"
$ws.Cells.Item(349, 10).Value = 1350
$ws.Cells.Item(349, 11).Value = 7000
$ws.Cells.Item(349, 12).Value = 8000
$ws.Cells.Item(349, 13).Value = 7615
$ws.Cells.Item(349, 14).Value = "`$/docena de matas"
$ws.Cells.Item(349, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(349, 16).Value = 1269
$ws.Cells.Item(349, 17).Value = 6
$ws.Cells.Item(349, 18).Value = "Hortaliza"

# Row 350: same new observation date, "Segunda" quality.
$ws.Cells.Item(350, 1).Value = 6
$ws.Cells.Item(350, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(350, 3).Value = "Metropolitana"
$ws.Cells.Item(350, 4).Value = 44461
$ws.Cells.Item(350, 5).Value = 13
$ws.Cells.Item(350, 6).Value = 100112017
$ws.Cells.Item(350, 7).Value = "Apio"
$ws.Cells.Item(350, 8).Value = "Americana (o)"
$ws.Cells.Item(350, 9).Value = "Segunda"
$ws.Cells.Item(350, 10).Value = 650
$ws.Cells.Item(350, 11).Value = 5000
$ws.Cells.Item(350, 12).Value = 5000
$ws.Cells.Item(350, 13).Value = 5000
$ws.Cells.Item(350, 14).Value = "`$/docena de matas"
$ws.Cells.Item(350, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(350, 16).Value = 833
$ws.Cells.Item(350, 17).Value = 6
$ws.Cells.Item(350, 18).Value = "Hortaliza"
